# Complementos.xlsx edit: re-sort the "Tabla2" table (range F2:G57 on sheet
# "Hoja1") alphabetically (A->Z) on the "Cargos" column (G), which is what
# the author did in Excel. Row numbers in column F ("No") keep the sorted
# row's original value, and cell formatting travels with its row, exactly
# like Excel's native Table-sort behaviour.
#
# Also updates the view state left behind by that interactive session:
# the active cell/selection moves to F15 and the sheet is no longer
# scrolled down to show row 40 at the top.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Sort Tabla2 (F2:G57) by column G, ascending, header row included ---
$lo = $ws.ListObjects.Item("Tabla2")
$lo.Sort.SortFields.Clear()
$sortColumn = $lo.ListColumns.Item("Cargos").Range
$lo.Sort.SortFields.Add($sortColumn, 0, 1, 0, 0)
$lo.Sort.Header = 1
$lo.Sort.Apply()

# --- Restore the view: scroll back to the top and move the selection ---
$win = $wb.Windows.Item(1)
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("F15").Select()
